$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row of data appended after the existing profit-allocation history
# (run on 2025-10-08). Force column A to Text so the "MM/DD/YYYY" string
# is stored literally instead of being auto-converted to a date serial,
# matching the existing rows in the sheet.
$ws.Range("A37").NumberFormat = "@"
$ws.Range("A37").Value = "10/08/2025"
$ws.Range("B37").Value = 0.1528034018842741
$ws.Range("C37").Value = 0.8471965981157259
